# 7.8 History Card & Advanced Story
# Updates Chen's dialogue lines in the sheet and refreshes the view selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update dialogue / question text in column B ---
# New lines (appended to the shared-string table in this order)
$ws.Range("B14").Value = " <color=#00CC00>(I recall they did indeed arrive at the banquet hall at the same time.)</color>"
$ws.Range("B2").Value = "Sir,  I’ll leave it to you to find the culprit."
$ws.Range("B5").Value = "When was the last time you saw the Lordr?"
$ws.Range("B6").Value = "Sometime after 5 PM. I was walking down the corridor and happened to see the master coming out of the main hall."
$ws.Range("B8").Value = "He often guided me in martial arts before."
$ws.Range("B10").Value = "What did you do before and after the banquet started?"
$ws.Range("B12").Value = " <color=#00CC00>(Butler He and Ling both mentioned this in their statements.)</color>"
$ws.Range("B16").Value = "Does going to the restroom in the banquet hall count? I was gone for about 15 min."

# Lines that re-use text already present elsewhere in the sheet
$ws.Range("B3").Value = "Very well. To uncover the truth, I need to ask you a few questions."
$ws.Range("B4").Value = "Ask me anything!"
$ws.Range("B7").Value = "I greeted him, then headed to the training ground to practice martial arts."
$ws.Range("B9").Value = "I never imagined that would be the last time I saw him."
$ws.Range("B11").Value = "After my training, I went back to my room to change into a fresh set of clothes, then went to the kitchen to help Ling prepare the food."
$ws.Range("B13").Value = "After cooking, I had just stepped out of the kitchen when I ran into Quan. We headed to the banquet hall together."
$ws.Range("B15").Value = "Did you leave the banquet at any point?"
$ws.Range("B17").Value = "Aside from that, I didn’t leave at all."

# --- Row 5 is now a single short line instead of a wrapped one ---
$ws.Rows.Item(5).RowHeight = 17

# --- Update the active view / selection ---
$win = $excel.ActiveWindow
$ws.Range("B23").Select()
$win.ScrollRow = 6
$win.ScrollColumn = 1
